# Mark meetreq_status (K) and fdh_status (L) as "done" for all data rows
# (rows 2-12) on the "all" sheet. Fixes bug where wrong time/room meant
# the meeting-room booking step never got flagged as completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 11).Value = "done"   # column K - meetreq_status
    $ws.Cells.Item($row, 12).Value = "done"   # column L - fdh_status
}
